# Adds the "conceptPath" values (column J) for each hospital-variable row,
# matching the new sharedStrings entries introduced by the commit
# "created an improved view for hospital variables".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value  = "/root/pet/av45"
$ws.Range("J3").Value  = "/root/pet/fdg"
$ws.Range("J4").Value  = "/root/pet/pib"
$ws.Range("J5").Value  = "/root/brain_anatomy/brainstem"
$ws.Range("J6").Value  = "/root/brain_anatomy/tiv"
$ws.Range("J7").Value  = "/root/brain_anatomy/csf_volume/_3rdventricle"
$ws.Range("J8").Value  = "/root/brain_anatomy/csf_volume/_4thventricle"
$ws.Range("J9").Value  = "/root/brain_anatomy/csf_volume/csfglobal"
$ws.Range("J10").Value = "/root/brain_anatomy/csf_volume/leftinflatvent"
$ws.Range("J11").Value = "/root/brain_anatomy/csf_volume/leftlateralventricle"
$ws.Range("J12").Value = "/root/brain_anatomy/csf_volume/rightinflatvent"
$ws.Range("J13").Value = "/root/brain_anatomy/csf_volume/rightlateralventricle"
$ws.Range("J14").Value = "/root/brain_anatomy/grey_matter_volume/cerebellum/cerebellarvermallobulesiv"
$ws.Range("J15").Value = "/root/brain_anatomy/grey_matter_volume/cerebellum/cerebellarvermallobulesviiix"
$ws.Range("J16").Value = "/root/brain_anatomy/grey_matter_volume/cerebellum/cerebellarvermallobulesvivii"
$ws.Range("J17").Value = "/root/brain_anatomy/grey_matter_volume/cerebellum/leftcerebellumexterior"

# Slightly narrower default columns and a dedicated, wider column for "comments" (I),
# matching the refreshed view layout.
$ws.Range("A1").EntireColumn.ColumnWidth = 12.53
$ws.Range("B1").EntireColumn.ColumnWidth = 40.74
$ws.Range("C1:H1").EntireColumn.ColumnWidth = 12.53
$ws.Range("I1").EntireColumn.ColumnWidth = 28.05
$ws.Range("J1").EntireColumn.ColumnWidth = 16.45

# New cursor position left behind after the edits.
$ws.Range("H29").Select()

"done"
